$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.9170291623044325
$ws.Range("D2").Value = 0.9416429316346989
$ws.Range("C3").Value = 0.8484261885130215
$ws.Range("D3").Value = 0.8912106471290957
$ws.Range("C4").Value = 0.7911289445959027
$ws.Range("D4").Value = 0.8476216671120647
$ws.Range("C5").Value = 0.7429989570171585
$ws.Range("D5").Value = 0.8097123667562878
$ws.Range("C6").Value = 0.702138676865307
$ws.Range("D6").Value = 0.7765327088450696
$ws.Range("C7").Value = 0.6668797697886445
$ws.Range("D7").Value = 0.7470442273258214
$ws.Range("C8").Value = 0.6369154434875203
$ws.Range("D8").Value = 0.720767687345103
$ws.Range("C9").Value = 0.6112691090540544
$ws.Range("D9").Value = 0.6974311988184334
$ws.Range("C10").Value = 0.5890682688138595
$ws.Range("D10").Value = 0.6776920056548507
$ws.Range("C11").Value = 0.5693487441109321
$ws.Range("D11").Value = 0.6602303130431306
$ws.Range("C12").Value = 0.5517831478055472
$ws.Range("D12").Value = 0.6439505993913955
$ws.Range("C13").Value = 0.5362537887404351
$ws.Range("D13").Value = 0.6296360109948211
$ws.Range("C14").Value = 0.5225059381354843
$ws.Range("D14").Value = 0.6170961508821817
$ws.Range("C15").Value = 0.5102461613608638
$ws.Range("D15").Value = 0.6055784430817203
$ws.Range("C16").Value = 0.4998193176503996
$ws.Range("D16").Value = 0.5952556561342357
$ws.Range("C17").Value = 0.4901817193478339
$ws.Range("D17").Value = 0.5857725507238485
$ws.Range("C18").Value = 0.4822907789146031
$ws.Range("D18").Value = 0.5776090876804424
$ws.Range("C19").Value = 0.474630959531397
$ws.Range("D19").Value = 0.5700648111644333
$ws.Range("C20").Value = 0.4685077653780572
$ws.Range("D20").Value = 0.5632518243206505
$ws.Range("C21").Value = 0.4627755445589086
$ws.Range("D21").Value = 0.5572483003596687
$ws.Range("C22").Value = 0.4578504594884943
$ws.Range("D22").Value = 0.5516873748309649
$ws.Range("C23").Value = 0.4533645636361425
$ws.Range("D23").Value = 0.5467834880525162
$ws.Range("C24").Value = 0.449385863586327
$ws.Range("D24").Value = 0.5431425355425593
$ws.Range("C25").Value = 0.4455361928658126
$ws.Range("D25").Value = 0.5392074880422915
$ws.Range("C26").Value = 0.4387220851264734
$ws.Range("D26").Value = 0.5326134768152555
$ws.Range("C27").Value = 0.4361911910186426
$ws.Range("D27").Value = 0.5302161287583528
$ws.Range("C28").Value = 0.4335698383388157
$ws.Range("D28").Value = 0.528035818947681
$ws.Range("C29").Value = 0.431149493616141
$ws.Range("D29").Value = 0.5256102068656964
$ws.Range("C30").Value = 0.4289546438603017
$ws.Range("D30").Value = 0.5231445996552525
$ws.Range("C31").Value = 0.4269719063580684
$ws.Range("D31").Value = 0.5213287972853136
$ws.Range("C32").Value = 0.4219043152845512
$ws.Range("D32").Value = 0.5168468522990183
$ws.Range("C33").Value = 0.4159431785730546
$ws.Range("D33").Value = 0.512201145119209
$ws.Range("C34").Value = 0.4152095396326018
$ws.Range("D34").Value = 0.5114757846556557
$ws.Range("C35").Value = 0.4145690827668099
$ws.Range("D35").Value = 0.5106911068106462
$ws.Range("C36").Value = 0.413990544795924
$ws.Range("D36").Value = 0.5100275917815962
$ws.Range("C37").Value = 0.410696331960496
$ws.Range("D37").Value = 0.5067050919967161
$ws.Range("C38").Value = 0.409717577510233
$ws.Range("D38").Value = 0.5057404340940624
$ws.Range("C39").Value = 0.4071357892940921
$ws.Range("D39").Value = 0.5037252381921539
